$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 446.92307
$ws.Range("I11").Value = 446.92307
$ws.Range("K11").Value = 446.92307
$ws.Range("M11").Value = -306.92307

$ws.Range("H12").Value = 907.8333
$ws.Range("I12").Value = 727.7143
$ws.Range("J12").Value = 1160
$ws.Range("K12").Value = 727.7143
$ws.Range("L12").Value = 1160
$ws.Range("M12").Value = -557.7143
$ws.Range("N12").Value = -1500

$ws.Range("H28").Value = 35714780
$ws.Range("I28").Value = 55555852
$ws.Range("J28").Value = 852.1
$ws.Range("K28").Value = 55555852
$ws.Range("L28").Value = 852.1
$ws.Range("M28").Value = -55555367
$ws.Range("N28").Value = -1822.1

$ws.Range("H57").Value = 135180.14
$ws.Range("J57").Value = 135180.14
$ws.Range("L57").Value = 405540.42
$ws.Range("N57").Value = -406538.42

$ws.Range("H64").Value = 4316.25
$ws.Range("I64").Value = 3369.8333
$ws.Range("J64").Value = 5262.6665
$ws.Range("K64").Value = 3369.8333
$ws.Range("L64").Value = 5262.6665
$ws.Range("M64").Value = -3121.8333
$ws.Range("N64").Value = -5758.6665

$ws.Range("H67").Value = 4316.25
$ws.Range("I67").Value = 3369.8333
$ws.Range("J67").Value = 5262.6665
$ws.Range("K67").Value = 3369.8333
$ws.Range("L67").Value = 5262.6665
$ws.Range("M67").Value = -2511.8333
$ws.Range("N67").Value = -6978.6665

$ws.Range("H93").Value = 28398.334
$ws.Range("J93").Value = 28398.334
$ws.Range("L93").Value = 28398.334
$ws.Range("N93").Value = -33390.334

$ws.Range("H95").Value = 29217.5
$ws.Range("J95").Value = 29217.5
$ws.Range("L95").Value = 29217.5
$ws.Range("N95").Value = -34709.5

$ws.Range("H107").Value = 14714322
$ws.Range("I107").Value = 20843908
$ws.Range("J107").Value = 3317.8
$ws.Range("K107").Value = 20843908
$ws.Range("L107").Value = 3317.8
$ws.Range("M107").Value = -20841988
$ws.Range("N107").Value = -7157.8

$ws.Range("H116").Value = 10044.444
$ws.Range("J116").Value = 10687.5
$ws.Range("L116").Value = 10687.5
$ws.Range("N116").Value = -17571.5

$ws.Range("H129").Value = 599.6667
$ws.Range("I129").Value = 599.6667
$ws.Range("K129").Value = 1799.0001
$ws.Range("M129").Value = 3200.9999

$ws.Range("H137").Value = 1286.3462
$ws.Range("I137").Value = 1313.4667
$ws.Range("J137").Value = 1249.3636
$ws.Range("K137").Value = 3940.4001
$ws.Range("L137").Value = 3748.0908
$ws.Range("M137").Value = -1390.4001
$ws.Range("N137").Value = -8848.0908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3722.052
$ws.Range("I32").Value = 2348.352
$ws.Range("K32").Value = 2348.352
$ws.Range("M32").Value = -2061.352

$ws.Range("H74").Value = 7855.407
$ws.Range("I74").Value = 9025.5
$ws.Range("J74").Value = 4512.2856
$ws.Range("K74").Value = 9025.5
$ws.Range("L74").Value = 4512.2856
$ws.Range("M74").Value = -8151.5
$ws.Range("N74").Value = -6260.2856

$ws.Range("H77").Value = 7855.407
$ws.Range("I77").Value = 9025.5
$ws.Range("J77").Value = 4512.2856
$ws.Range("K77").Value = 45127.5
$ws.Range("L77").Value = 22561.428
$ws.Range("M77").Value = -40759.5
$ws.Range("N77").Value = -31297.428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 9463.5
$ws.Range("I96").Value = 9463.5
$ws.Range("K96").Value = 9463.5
$ws.Range("M96").Value = -6717.5

$ws.Range("H134").Value = 6823.7666
$ws.Range("I134").Value = 6474.6665
$ws.Range("J134").Value = 9965.666999999999
$ws.Range("K134").Value = 19423.9995
$ws.Range("L134").Value = 29897.001
$ws.Range("M134").Value = -16888.9995
$ws.Range("N134").Value = -34967.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7084.2915
$ws.Range("I58").Value = 6219.353
$ws.Range("J58").Value = 9184.857
$ws.Range("K58").Value = 6219.353
$ws.Range("L58").Value = 9184.857
$ws.Range("M58").Value = -6016.353
$ws.Range("N58").Value = -9590.857

$ws.Range("H132").Value = 11045
$ws.Range("I132").Value = 13371.8
$ws.Range("J132").Value = 4397
$ws.Range("K132").Value = 40115.39999999999
$ws.Range("L132").Value = 13191
$ws.Range("M132").Value = -37585.39999999999
$ws.Range("N132").Value = -18251

$ws.Range("H134").Value = 10448.27
$ws.Range("I134").Value = 14083.0625
$ws.Range("J134").Value = 4632.6
$ws.Range("K134").Value = 42249.1875
$ws.Range("L134").Value = 13897.8
$ws.Range("M134").Value = -39714.1875
$ws.Range("N134").Value = -18967.8

$ws.Range("H136").Value = 7084.2915
$ws.Range("I136").Value = 6219.353
$ws.Range("J136").Value = 9184.857
$ws.Range("K136").Value = 18658.059
$ws.Range("L136").Value = 27554.571
$ws.Range("M136").Value = -16108.059
$ws.Range("N136").Value = -32654.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1059.375
$ws.Range("J5").Value = 2500
$ws.Range("L5").Value = 7500
$ws.Range("N5").Value = -7724

$ws.Range("H19").Value = 400
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H114").Value = 2318.3333
$ws.Range("I114").Value = 2000
$ws.Range("J114").Value = 2477.5
$ws.Range("K114").Value = 6000
$ws.Range("L114").Value = 7432.5
$ws.Range("M114").Value = -2746
$ws.Range("N114").Value = -13940.5

$ws.Range("H117").Value = 2504.3635
$ws.Range("I117").Value = 1110.8
$ws.Range("J117").Value = 3665.6667
$ws.Range("K117").Value = 3332.4
$ws.Range("L117").Value = 10997.0001
$ws.Range("M117").Value = 109.6000000000004
$ws.Range("N117").Value = -17881.0001

$ws.Range("H135").Value = 1059.375
$ws.Range("J135").Value = 2500
$ws.Range("L135").Value = 22500
$ws.Range("N135").Value = -27570

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 3397.5715
$ws.Range("I31").Value = 4350
$ws.Range("J31").Value = 3016.6
$ws.Range("K31").Value = 4350
$ws.Range("L31").Value = 3016.6
$ws.Range("M31").Value = -4102
$ws.Range("N31").Value = -3512.6

$ws.Range("H40").Value = 15216.538
$ws.Range("I40").Value = 16290.789
$ws.Range("K40").Value = 16290.789
$ws.Range("M40").Value = -16154.789

$ws.Range("H46").Value = 3418.0952

$ws.Range("H68").Value = 14496824
$ws.Range("I68").Value = 17545882
$ws.Range("K68").Value = 17545882
$ws.Range("M68").Value = -17545133

$ws.Range("H71").Value = 14496824
$ws.Range("I71").Value = 17545882
$ws.Range("K71").Value = 87729410
$ws.Range("M71").Value = -87725666

$ws.Range("H93").Value = 9525060
$ws.Range("I93").Value = 11112236
$ws.Range("K93").Value = 11112236
$ws.Range("M93").Value = -11110988

$ws.Range("H122").Value = 7940.9697
$ws.Range("I122").Value = 7931.375
$ws.Range("K122").Value = 23794.125
$ws.Range("M122").Value = -21344.125

$ws.Range("H130").Value = 179999
$ws.Range("J130").Value = 179999
$ws.Range("L130").Value = 179999
$ws.Range("N130").Value = -190039

$ws.Range("H136").Value = 6675327
$ws.Range("J136").Value = 15062.5
$ws.Range("L136").Value = 45187.5
$ws.Range("N136").Value = -50287.5

$ws.Range("H141").Value = 150000
$ws.Range("J141").Value = 150000
$ws.Range("L141").Value = 150000
$ws.Range("N141").Value = -160360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 8756

$ws.Range("H113").Value = 528.75
$ws.Range("I113").Value = 345.375
$ws.Range("K113").Value = 1036.125
$ws.Range("M113").Value = 1133.875

$ws.Range("H135").Value = 218265
$ws.Range("J135").Value = 218265
$ws.Range("L135").Value = 218265
$ws.Range("N135").Value = -228405
